# Logged Week 16 and performed season sim from Week 17
$wb = $excel.ActiveWorkbook

$wsQB = $wb.Worksheets.Item("QB")
$wsWR = $wb.Worksheets.Item("WR")

# Week 16 stat logging added a new WR to the roster: D.Sills, with a
# fresh (all-zero) stat line appended to the bottom of the table.
$wsWR.Range("A10").Value = "D.Sills"
$wsWR.Range("B10:J10").Value = 0

# Leave the QB sheet's old selection as-is, but it is no longer the
# tab that is active/selected in the workbook.
$wsQB.Range("L5").Select()

# The season sim from Week 17 finished up on the WR sheet, which is
# now the active/selected tab, sitting just past the newly added row.
$wsWR.Activate()
$wsWR.Range("J11").Select()
